$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 180, shifting existing rows 180-183 down to 181-184.
$ws.Rows.Item(180).Insert()

# Fill in the new row 180 with the new weekly record.
$ws.Range("A180").Value = 8
$ws.Range("B180").Value = "Terminal La Palmera de La Serena"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 44448
$ws.Range("E180").Value = 4
$ws.Range("F180").Value = 100114013
$ws.Range("G180").Value = "Zanahoria"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 720
$ws.Range("K180").Value = 4500
$ws.Range("L180").Value = 5000
$ws.Range("M180").Value = 4750
$ws.Range("N180").Value = "`$/saco 20 kilos"
$ws.Range("O180").Value = "Provincia del Elquí"
$ws.Range("P180").Value = 238
$ws.Range("Q180").Value = 20
$ws.Range("R180").Value = "Hortaliza"
